$wb = $excel.ActiveWorkbook

# Sheet ALC, row 107 (anchor G=27766)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 970
$ws.Range("I107").Value = 871.1429000000001
$ws.Range("J107").Value = 1068.8572
$ws.Range("K107").Value = 871.1429000000001
$ws.Range("L107").Value = 1068.8572
$ws.Range("M107").Value = 1048.8571
$ws.Range("N107").Value = -4908.8572

# Sheet ALC, row 137 (anchor G=44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2063.923
$ws.Range("I137").Value = 1436.7778
$ws.Range("J137").Value = 3475
$ws.Range("K137").Value = 4310.3334
$ws.Range("L137").Value = 10425
$ws.Range("M137").Value = -1760.3334
$ws.Range("N137").Value = -15525

# Sheet ARM, row 32 (anchor G=44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16454.744
$ws.Range("I32").Value = 17719.03
$ws.Range("J32").Value = 9079.75
$ws.Range("K32").Value = 17719.03
$ws.Range("L32").Value = 9079.75
$ws.Range("M32").Value = -17432.03
$ws.Range("N32").Value = -9653.75

# Sheet ARM, row 61 (anchor G=43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7540.3726
$ws.Range("I61").Value = 5768.8613
$ws.Range("K61").Value = 5768.8613
$ws.Range("M61").Value = -5556.8613

# Sheet ARM, row 122 (anchor G=36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4808940
$ws.Range("I122").Value = 1235.1052
$ws.Range("J122").Value = 17858424
$ws.Range("K122").Value = 3705.3156
$ws.Range("L122").Value = 53575272
$ws.Range("M122").Value = -1255.3156
$ws.Range("N122").Value = -53580172

# Sheet ARM, row 132 (anchor G=43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 8381.714
$ws.Range("I132").Value = 3641.7778
$ws.Range("J132").Value = 11936.667
$ws.Range("K132").Value = 10925.3334
$ws.Range("L132").Value = 35810.001
$ws.Range("M132").Value = -8395.3334
$ws.Range("N132").Value = -40870.001

# Sheet ARM, row 136 (anchor G=43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7540.3726
$ws.Range("I136").Value = 5768.8613
$ws.Range("K136").Value = 17306.5839
$ws.Range("M136").Value = -14756.5839

# Sheet BSM, row 107 (anchor G=27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1509.5294
$ws.Range("I107").Value = 936.2
$ws.Range("J107").Value = 2328.5715
$ws.Range("K107").Value = 936.2
$ws.Range("L107").Value = 2328.5715
$ws.Range("M107").Value = 983.8
$ws.Range("N107").Value = -6168.5715

# Sheet CRP, row 16 (anchor G=27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1339.4615
$ws.Range("I16").Value = 666.6667
$ws.Range("J16").Value = 1916.1428
$ws.Range("K16").Value = 666.6667
$ws.Range("L16").Value = 1916.1428
$ws.Range("M16").Value = -379.6667
$ws.Range("N16").Value = -2490.1428

# Sheet CRP, row 22 (anchor G=5367)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 246.83333
$ws.Range("I22").Value = 200.5
$ws.Range("J22").Value = 270
$ws.Range("K22").Value = 200.5
$ws.Range("L22").Value = 270
$ws.Range("M22").Value = 149.5
$ws.Range("N22").Value = -970

# Sheet CRP, row 58 (anchor G=44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1422470.4
$ws.Range("I58").Value = 2165899.5
$ws.Range("J58").Value = 3196.818
$ws.Range("K58").Value = 2165899.5
$ws.Range("L58").Value = 3196.818
$ws.Range("M58").Value = -2165696.5
$ws.Range("N58").Value = -3602.818

# Sheet CRP, row 113 (anchor G=27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1339.4615
$ws.Range("I113").Value = 666.6667
$ws.Range("J113").Value = 1916.1428
$ws.Range("K113").Value = 666.6667
$ws.Range("L113").Value = 1916.1428
$ws.Range("M113").Value = 1503.3333
$ws.Range("N113").Value = -6256.1428

# Sheet CRP, row 132 (anchor G=44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7436
$ws.Range("I132").Value = 10300.154
$ws.Range("K132").Value = 30900.462
$ws.Range("M132").Value = -28370.462

# Sheet CRP, row 136 (anchor G=44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1422470.4
$ws.Range("I136").Value = 2165899.5
$ws.Range("J136").Value = 3196.818
$ws.Range("K136").Value = 6497698.5
$ws.Range("L136").Value = 9590.454000000002
$ws.Range("M136").Value = -6495148.5
$ws.Range("N136").Value = -14690.454

# Sheet CUL, row 114 (anchor G=27865)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 3704426.2
$ws.Range("I114").Value = 338.5
$ws.Range("J114").Value = 7937669.5
$ws.Range("K114").Value = 1015.5
$ws.Range("L114").Value = 23813008.5
$ws.Range("M114").Value = 2238.5
$ws.Range("N114").Value = -23819516.5

# Sheet CUL, row 117 (anchor G=27870)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 37037836
$ws.Range("I117").Value = 400
$ws.Range("J117").Value = 55556556
$ws.Range("K117").Value = 1200
$ws.Range("L117").Value = 166669668
$ws.Range("M117").Value = 2242
$ws.Range("N117").Value = -166676552

# Sheet CUL, row 121 (anchor G=27878)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1147.2858
$ws.Range("I121").Value = 330
$ws.Range("J121").Value = 1283.5
$ws.Range("K121").Value = 990
$ws.Range("L121").Value = 3850.5
$ws.Range("M121").Value = 320
$ws.Range("N121").Value = -6470.5

# Sheet CUL, row 122 (anchor G=36078)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 830.65625
$ws.Range("I122").Value = 650.9091
$ws.Range("J122").Value = 924.8095
$ws.Range("K122").Value = 5858.1819
$ws.Range("L122").Value = 8323.2855
$ws.Range("M122").Value = -3408.1819
$ws.Range("N122").Value = -13223.2855

# Sheet CUL, row 131 (anchor G=36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 26027.691
$ws.Range("J131").Value = 29378.234
$ws.Range("L131").Value = 88134.702
$ws.Range("N131").Value = -98214.702

# Sheet CUL, row 139 (anchor G=44102)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1533398.1
$ws.Range("I139").Value = 2711017.5
$ws.Range("J139").Value = 2493.2
$ws.Range("K139").Value = 8133052.5
$ws.Range("L139").Value = 7479.599999999999
$ws.Range("M139").Value = -8127912.5
$ws.Range("N139").Value = -17759.6

# Sheet GSM, row 53 (anchor G=4361)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 14900
$ws.Range("J53").Value = 14900
$ws.Range("L53").Value = 14900
$ws.Range("M53").Value = -16162

# Sheet GSM, row 126 (anchor G=36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2561.375
$ws.Range("I126").Value = 1634.6154
$ws.Range("J126").Value = 3195.4736
$ws.Range("K126").Value = 4903.8462
$ws.Range("L126").Value = 9586.4208
$ws.Range("M126").Value = -2433.8462
$ws.Range("N126").Value = -14526.4208

# Sheet GSM, row 132 (anchor G=44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6845
$ws.Range("I132").Value = 3392.4443
$ws.Range("J132").Value = 14613.25
$ws.Range("K132").Value = 10177.3329
$ws.Range("L132").Value = 43839.75
$ws.Range("M132").Value = -7647.332900000001
$ws.Range("N132").Value = -48899.75

# Sheet GSM, row 133 (anchor G=41854)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 85000
$ws.Range("J133").Value = 85000
$ws.Range("L133").Value = 85000
$ws.Range("N133").Value = -95120

# Sheet LTW, row 132 (anchor G=44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5199.6665
$ws.Range("I132").Value = 5315.577
$ws.Range("J132").Value = 4898.3
$ws.Range("K132").Value = 15946.731
$ws.Range("L132").Value = 14694.9
$ws.Range("M132").Value = -13416.731
$ws.Range("N132").Value = -19754.9

# Sheet WVR, row 7 (anchor G=2661)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 21400
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 21400
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 21400
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -21626

# Sheet WVR, row 12 (anchor G=3316)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 2511025
$ws.Range("J12").Value = 14700
$ws.Range("L12").Value = 14700
$ws.Range("N12").Value = -14984

# Sheet WVR, row 132 (anchor G=44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2956.75
$ws.Range("I132").Value = 2941.6584
$ws.Range("J132").Value = 2989.3157
$ws.Range("K132").Value = 8824.975199999999
$ws.Range("L132").Value = 8967.947100000001
$ws.Range("M132").Value = -6294.975199999999
$ws.Range("N132").Value = -14027.9471

# Sheet WVR, row 136 (anchor G=44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4261.5835
$ws.Range("I136").Value = 3192.9246
$ws.Range("J136").Value = 6088.645
$ws.Range("K136").Value = 9578.773799999999
$ws.Range("L136").Value = 18265.935
$ws.Range("M136").Value = -7028.773799999999
$ws.Range("N136").Value = -23365.935
